# Price/Hora refresh + coin-list reshuffle for the 22-12-2022 13:xx snapshot.
# (Row 10 gains a new top entry - 'One' - pushing WazirX..CoinExToken down by one row;
#  every row's 'Hora' (G) ticks from 12 -> 13; several 'Price' (D) values are refreshed.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells that hold numeric-looking text (Price, Hora) must stay TEXT, matching the
# source sheet's inline-string cells - so we type them with a leading apostrophe
# (Excel's own 'force text' convention) and then strip the resulting cell style back
# to Normal so no stray formatting is left behind.
function Set-TextValue($range, [string]$value) {
    $range.Value = "'" + $value
    $range.Style = "Normal"
}

$rowUpdates = @(
    @{ Row=2; D="244.50"; G="13" }
    @{ Row=3; D="22.50"; G="13" }
    @{ Row=4; D="5.405"; G="13" }
    @{ Row=5; D="0.05775"; G="13" }
    @{ Row=6; D="3.430"; G="13" }
    @{ Row=7; D="6.327"; G="13" }
    @{ Row=8; D="0.8108"; G="13" }
    @{ Row=9; D="0.8827"; G="13" }
    @{ Row=10; B="One"; C="https://coinranking.com/coin/6Lga5NiXX3rT+one-one"; D="0.0005842"; E="9OneONE"; G="13" }
    @{ Row=11; B="WazirX"; C="https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"; D="0.1449"; E="10WazirXWRX"; G="13" }
    @{ Row=12; B="MandalaExchangeToken"; C="https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"; D="0.07339"; E="11MandalaExchangeTokenMDX"; G="13" }
    @{ Row=13; B="LiechtensteinCryptoassetsExchange"; C="https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"; D="0.03036"; E="12LiechtensteinCryptoassetsExchangeLCX"; G="13" }
    @{ Row=14; B="BitrueCoin"; C="https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"; D="0.03099"; E="13BitrueCoinBTR"; G="13" }
    @{ Row=15; B="BitMartToken"; C="https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"; D="0.09401"; E="14BitMartTokenBMX"; G="13" }
    @{ Row=16; B="BitForexToken"; C="https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"; D="0.001592"; E="15BitForexTokenBF"; G="13" }
    @{ Row=17; B="CoinExToken"; C="https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"; D="0.04838"; E="16CoinExTokenCET"; G="13" }
    @{ Row=18; D="0.006394"; G="13" }
    @{ Row=19; D="0.004137"; G="13" }
    @{ Row=20; D="0.0009954"; G="13" }
    @{ Row=21; D="0.0001500"; G="13" }
    @{ Row=22; G="13" }
    @{ Row=23; D="2.195"; G="13" }
    @{ Row=24; D="0.3277"; G="13" }
    @{ Row=25; D="0.1321"; G="13" }
    @{ Row=26; D="4.180"; G="13" }
    @{ Row=27; D="0.0004651"; G="13" }
    @{ Row=28; G="13" }
    @{ Row=29; G="13" }
    @{ Row=30; G="13" }
    @{ Row=31; G="13" }
    @{ Row=32; G="13" }
    @{ Row=33; G="13" }
    @{ Row=34; G="13" }
    @{ Row=35; G="13" }
    @{ Row=36; G="13" }
    @{ Row=37; G="13" }
    @{ Row=38; G="13" }
    @{ Row=39; G="13" }
    @{ Row=40; D="0.03899"; G="13" }
    @{ Row=41; G="13" }
    @{ Row=42; D="0.1069"; G="13" }
    @{ Row=43; D="0.002601"; G="13" }
    @{ Row=44; D="0.007318"; G="13" }
    @{ Row=45; D="0.00005592"; G="13" }
    @{ Row=46; G="13" }
    @{ Row=47; D="0.3801"; G="13" }
    @{ Row=48; D="0.1583"; G="13" }
    @{ Row=49; G="13" }
    @{ Row=50; D="0.01010"; G="13" }
    @{ Row=51; G="13" }
)

foreach ($u in $rowUpdates) {
    $r = $u.Row
    if ($u.ContainsKey("B")) { $ws.Range("B$r").Value = $u.B }
    if ($u.ContainsKey("C")) { $ws.Range("C$r").Value = $u.C }
    if ($u.ContainsKey("D")) { Set-TextValue $ws.Range("D$r") $u.D }
    if ($u.ContainsKey("E")) { $ws.Range("E$r").Value = $u.E }
    if ($u.ContainsKey("G")) { Set-TextValue $ws.Range("G$r") $u.G }
}

